$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.316.06'
$ws.Range("E2").Value = '  +0.24%  '

$ws.Range("D3").Value = '3.727.44'
$ws.Range("E3").Value = '  -2.12%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '599.26'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.67%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.81'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.16%  '

$ws.Range("D7").Value = '3.725.36'
$ws.Range("E7").Value = '  -2.13%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("E9").Value = '  +0.11%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.170'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.31%  '

$ws.Range("E11").Value = '  -1.35%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.463'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.66%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.50'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.23%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000248'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.43%  '

$ws.Range("D15").Value = '4.348.31'
$ws.Range("E15").Value = '  -2.00%  '

$ws.Range("D16").Value = '3.722.93'
$ws.Range("E16").Value = '  -1.78%  '

$ws.Range("D17").Value = '68.249.43'
$ws.Range("E17").Value = '  -0.09%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.35'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.31%  '

$ws.Range("E19").Value = '  -0.81%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.30'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +7.81%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '491.58'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.37%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.31'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.57%  '

$ws.Range("E23").Value = '  -1.38%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.04'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.08%  '

$ws.Range("B25").Value = 'PEPE'
$ws.Range("C25").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000144'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.66%  '

$ws.Range("B26").Value = 'Fetch.AI'
$ws.Range("C26").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.33'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.61%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.35'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.63%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.16'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.23%  '

$ws.Range("E29").Value = '  -0.03%  '

$ws.Range("E30").Value = '  -0.02%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.89'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.16%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.39'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.14%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.65'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.17%  '

$ws.Range("D34").Value = '3.868.01'
$ws.Range("E34").Value = '  -1.90%  '

$ws.Range("E35").Value = '  -1.65%  '

$ws.Range("D36").Value = '3.671.20'
$ws.Range("E36").Value = '  -2.01%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.12%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.06%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.86'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.52%  '

$ws.Range("E40").Value = '  -2.43%  '

$ws.Range("E41").Value = '  -0.25%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '433.72'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.39%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '48.71'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.96%  '

$ws.Range("E44").Value = '  -1.77%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.87'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.86%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.47'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.48%  '

$ws.Range("B47").Value = 'USDe'
$ws.Range("C47").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.00'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.02%  '

$ws.Range("B48").Value = 'Arweave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '40.62'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.07%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.29'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.56%  '

$ws.Range("E50").Value = '  +0.30%  '

$ws.Range("D51").Value = '2.761.12'
$ws.Range("E51").Value = '  -2.88%  '
